$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.952.42"
$ws.Range("E2").Value = "  -0.08%  "

# Row 3
$ws.Range("D3").Value = "1.639.40"
$ws.Range("E3").Value = "  -0.16%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.46"
$ws.Range("E5").Value = "  -0.04%  "

# Row 6
$ws.Range("E6").Value = "  -0.43%  "

# Row 7
$ws.Range("E7").Value = "  -0.06%  "

# Row 8
$ws.Range("E8").Value = "  -0.74%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.258"
$ws.Range("E9").Value = "  -2.39%  "

# Row 10
$ws.Range("E10").Value = "  +0.00%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0885"
$ws.Range("E11").Value = "  +1.69%  "

# Row 12
$ws.Range("D12").Value = "1.871.84"

# Row 13
$ws.Range("D13").Value = "1.639.18"
$ws.Range("E13").Value = "  -0.37%  "

# Row 14
$ws.Range("E14").Value = "  +0.18%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.570"
$ws.Range("E15").Value = "  +0.89%  "

# Row 16
$ws.Range("E16").Value = "  -0.27%  "

# Row 17
$ws.Range("D17").Value = "27.950.88"
$ws.Range("E17").Value = "  +0.00%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "233.06"
$ws.Range("E18").Value = "  +0.45%  "

# Row 19
$ws.Range("E19").Value = "  -0.26%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.55"
$ws.Range("E20").Value = "  -1.46%  "

# Row 21
$ws.Range("E21").Value = "  +0.03%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.47"
$ws.Range("E22").Value = "  -2.69%  "

# Row 23
$ws.Range("E23").Value = "  -0.69%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.08"
$ws.Range("E24").Value = "  -3.57%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.03"
$ws.Range("E25").Value = "  +1.26%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.91"
$ws.Range("E26").Value = "  -0.19%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.65"
$ws.Range("E27").Value = "  -0.37%  "

# Row 28
$ws.Range("E28").Value = "  -0.69%  "

# Row 29
$ws.Range("E29").Value = "  -0.05%  "

# Row 30
$ws.Range("E30").Value = "  +0.00%  "

# Row 31
$ws.Range("E31").Value = "  +0.25%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.41"
$ws.Range("E32").Value = "  +2.95%  "

# Row 33
$ws.Range("D33").Value = "1.404.99"
$ws.Range("E33").Value = "  -4.06%  "

# Row 34
$ws.Range("E34").Value = "  -0.33%  "

# Row 35
$ws.Range("E35").Value = "  +1.50%  "

# Row 36
$ws.Range("E36").Value = "  +1.29%  "

# Row 37
$ws.Range("E37").Value = "  +0.37%  "

# Row 38
$ws.Range("E38").Value = "  +0.17%  "

# Row 39
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.927"
$ws.Range("E39").Value = "  -0.51%  "

# Row 40
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.878"
$ws.Range("E40").Value = "  -1.49%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.02"
$ws.Range("E41").Value = "  +0.73%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "67.08"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.85"
$ws.Range("E44").Value = "  +2.20%  "

# Row 45
$ws.Range("E45").Value = "  +2.91%  "

# Row 46
$ws.Range("E46").Value = "  -0.53%  "

# Row 47
$ws.Range("D47").Value = "1.780.84"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.94"
$ws.Range("E48").Value = "  -0.17%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.100"
$ws.Range("E49").Value = "  -0.40%  "

# Row 50
$ws.Range("E50").Value = "  -0.12%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.56"
$ws.Range("E51").Value = "  -2.05%  "
